$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.839.14"
$ws.Range("E2").Value = "  +5.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.897.04"
$ws.Range("E3").Value = "  +3.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.26"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4719"
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("E8").Value = "  +6.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.83"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08132"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.33"
$ws.Range("E12").Value = "  +6.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.078"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.870.62"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.330"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.21"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001051"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06633"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.872.68"
$ws.Range("E22").Value = "  +5.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.541"
$ws.Range("E23").Value = "  +3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.100.32"
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.65"
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.147"
$ws.Range("E29").Value = "  +4.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.534"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.52"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.006"
$ws.Range("E32").Value = "  +6.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09564"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.658"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("E35").Value = "  +6.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.396"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06193"
$ws.Range("E37").Value = "  +4.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02282"
$ws.Range("E38").Value = "  +4.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.608"
$ws.Range("E39").Value = "  +7.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.198"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5996"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1895"
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.45"
$ws.Range("E44").Value = "  +4.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.282"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5613"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.21"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.975"
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07276"
$ws.Range("E49").Value = "  +10.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.112"
$ws.Range("E50").Value = "  +15.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.75"
$ws.Range("E51").Value = "  +1.95%  "
